# Auto-generated COM-interop script implementing the '2022-Q1' sheet addition.
$wb = $excel.ActiveWorkbook

# --- Step 1: the sheet currently named "总计" (6th tab) keeps its file/rId slot
# but is repurposed to hold the new quarter's fund-holdings table, and is
# renamed to "2022-Q1". Grab a reference to its pre-existing header/index
# style (style index 2 in the source workbook) before we touch any values,
# so later cells can inherit the exact same formatting via Copy/PasteSpecial.
$quarterSheet = $wb.Worksheets.Item(6)
$styleSrc = $quarterSheet.Range("D1")
$idxStyleSrc = $quarterSheet.Range("A2")
$quarterSheet.Name = '2022-Q1'

# --- Step 2: insert a brand-new sheet right after it, taking over the
# "总计" name and the roll-up table that used to live in the old sheet.
$summarySheet = $wb.Worksheets.Add($null, $quarterSheet)
$summarySheet.Name = '总计'

# ===================== 2022-Q1 fund holdings table =====================
$hdr = $quarterSheet
# extend the bold/bordered header style (already on B1:D1) across E1:H1
$styleSrc.Copy() | Out-Null
$hdr.Range("E1:H1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$hdr.Cells.Item(1,2).Value = '基金代码'
$hdr.Cells.Item(1,3).Value = '基金名称'
$hdr.Cells.Item(1,4).Value = '基金规模'
$hdr.Cells.Item(1,5).Value = '股票总仓位'
$hdr.Cells.Item(1,6).Value = '仓位占比'
$hdr.Cells.Item(1,7).Value = '持有市值(亿元)'
$hdr.Cells.Item(1,8).Value = '仓位排名'

# column-A row index cells (A2:A17) use the same centered/bordered style
# the original A2:A6 already carried; extend it down to the new rows first.
$idxStyleSrc.Copy() | Out-Null
$hdr.Range("A2:A17").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# row 2: 320011
$hdr.Cells.Item(2,1).Value = 0
$hdr.Cells.Item(2,2).Value = '''320011'
$hdr.Cells.Item(2,3).Value = '诺安中小盘精选混合'
$hdr.Cells.Item(2,4).Value = '''3.67'
$hdr.Cells.Item(2,5).Value = '''84.64'
$hdr.Cells.Item(2,6).Value = '''5.13'
$hdr.Cells.Item(2,7).Value = '''0.1883'
$hdr.Cells.Item(2,8).Value = 3
# row 3: 004350
$hdr.Cells.Item(3,1).Value = 1
$hdr.Cells.Item(3,2).Value = '''004350'
$hdr.Cells.Item(3,3).Value = '汇丰晋信价值先锋股票'
$hdr.Cells.Item(3,4).Value = '''4.99'
$hdr.Cells.Item(3,5).Value = '''93.32'
$hdr.Cells.Item(3,6).Value = '''3.72'
$hdr.Cells.Item(3,7).Value = '''0.1856'
$hdr.Cells.Item(3,8).Value = 2
# row 4: 510160
$hdr.Cells.Item(4,1).Value = 2
$hdr.Cells.Item(4,2).Value = '''510160'
$hdr.Cells.Item(4,3).Value = '南方中证南方小康产业ETF'
$hdr.Cells.Item(4,4).Value = '''2.60'
$hdr.Cells.Item(4,5).Value = '''99.43'
$hdr.Cells.Item(4,6).Value = '''6.80'
$hdr.Cells.Item(4,7).Value = '''0.1768'
$hdr.Cells.Item(4,8).Value = 1
# row 5: 159811
$hdr.Cells.Item(5,1).Value = 3
$hdr.Cells.Item(5,2).Value = '''159811'
$hdr.Cells.Item(5,3).Value = '博时中证5G产业50ETF'
$hdr.Cells.Item(5,4).Value = '''2.60'
$hdr.Cells.Item(5,5).Value = '''97.65'
$hdr.Cells.Item(5,6).Value = '''4.96'
$hdr.Cells.Item(5,7).Value = '''0.1290'
$hdr.Cells.Item(5,8).Value = 6
# row 6: 001528
$hdr.Cells.Item(6,1).Value = 4
$hdr.Cells.Item(6,2).Value = '''001528'
$hdr.Cells.Item(6,3).Value = '诺安先进制造股票'
$hdr.Cells.Item(6,4).Value = '''2.44'
$hdr.Cells.Item(6,5).Value = '''87.44'
$hdr.Cells.Item(6,6).Value = '''5.27'
$hdr.Cells.Item(6,7).Value = '''0.1286'
$hdr.Cells.Item(6,8).Value = 5
# row 7: 011132
$hdr.Cells.Item(7,1).Value = 5
$hdr.Cells.Item(7,2).Value = '''011132'
$hdr.Cells.Item(7,3).Value = '鹏扬沪深300质量成长低波动指数证券投资基金A'
$hdr.Cells.Item(7,4).Value = '''2.73'
$hdr.Cells.Item(7,5).Value = '''94.39'
$hdr.Cells.Item(7,6).Value = '''2.90'
$hdr.Cells.Item(7,7).Value = '''0.0792'
$hdr.Cells.Item(7,8).Value = 4
# row 8: 005870
$hdr.Cells.Item(8,1).Value = 6
$hdr.Cells.Item(8,2).Value = '''005870'
$hdr.Cells.Item(8,3).Value = '鹏华沪深300指数增强'
$hdr.Cells.Item(8,4).Value = '''3.53'
$hdr.Cells.Item(8,5).Value = '''92.67'
$hdr.Cells.Item(8,6).Value = '''2.21'
$hdr.Cells.Item(8,7).Value = '''0.0780'
$hdr.Cells.Item(8,8).Value = 7
# row 9: 320015
$hdr.Cells.Item(9,1).Value = 7
$hdr.Cells.Item(9,2).Value = '''320015'
$hdr.Cells.Item(9,3).Value = '诺安行业轮动混合'
$hdr.Cells.Item(9,4).Value = '''1.29'
$hdr.Cells.Item(9,5).Value = '''85.68'
$hdr.Cells.Item(9,6).Value = '''4.48'
$hdr.Cells.Item(9,7).Value = '''0.0578'
$hdr.Cells.Item(9,8).Value = 5
# row 10: 006429
$hdr.Cells.Item(10,1).Value = 8
$hdr.Cells.Item(10,2).Value = '''006429'
$hdr.Cells.Item(10,3).Value = '诺安恒鑫混合'
$hdr.Cells.Item(10,4).Value = '''0.82'
$hdr.Cells.Item(10,5).Value = '''85.51'
$hdr.Cells.Item(10,6).Value = '''5.51'
$hdr.Cells.Item(10,7).Value = '''0.0452'
$hdr.Cells.Item(10,8).Value = 3
# row 11: 011243
$hdr.Cells.Item(11,1).Value = 9
$hdr.Cells.Item(11,2).Value = '''011243'
$hdr.Cells.Item(11,3).Value = '万家惠裕回报6个月持有期混合型证券投资基金A'
$hdr.Cells.Item(11,4).Value = '''4.93'
$hdr.Cells.Item(11,5).Value = '''23.04'
$hdr.Cells.Item(11,6).Value = '''0.88'
$hdr.Cells.Item(11,7).Value = '''0.0434'
$hdr.Cells.Item(11,8).Value = 6
# row 12: 012377
$hdr.Cells.Item(12,1).Value = 10
$hdr.Cells.Item(12,2).Value = '''012377'
$hdr.Cells.Item(12,3).Value = '长盛安睿一年持有混合A'
$hdr.Cells.Item(12,4).Value = '''4.20'
$hdr.Cells.Item(12,5).Value = '''29.09'
$hdr.Cells.Item(12,6).Value = '''0.59'
$hdr.Cells.Item(12,7).Value = '''0.0248'
$hdr.Cells.Item(12,8).Value = 10
# row 13: 011133
$hdr.Cells.Item(13,1).Value = 11
$hdr.Cells.Item(13,2).Value = '''011133'
$hdr.Cells.Item(13,3).Value = '鹏扬沪深300质量成长低波动指数证券投资基金C'
$hdr.Cells.Item(13,4).Value = '''0.80'
$hdr.Cells.Item(13,5).Value = '''94.39'
$hdr.Cells.Item(13,6).Value = '''2.90'
$hdr.Cells.Item(13,7).Value = '''0.0232'
$hdr.Cells.Item(13,8).Value = 4
# row 14: 005035
$hdr.Cells.Item(14,1).Value = 12
$hdr.Cells.Item(14,2).Value = '''005035'
$hdr.Cells.Item(14,3).Value = '银华信息科技量化优选股票A'
$hdr.Cells.Item(14,4).Value = '''0.26'
$hdr.Cells.Item(14,5).Value = '''90.98'
$hdr.Cells.Item(14,6).Value = '''2.28'
$hdr.Cells.Item(14,7).Value = '''0.0059'
$hdr.Cells.Item(14,8).Value = 6
# row 15: 012378
$hdr.Cells.Item(15,1).Value = 13
$hdr.Cells.Item(15,2).Value = '''012378'
$hdr.Cells.Item(15,3).Value = '长盛安睿一年持有混合C'
$hdr.Cells.Item(15,4).Value = '''0.36'
$hdr.Cells.Item(15,5).Value = '''29.09'
$hdr.Cells.Item(15,6).Value = '''0.59'
$hdr.Cells.Item(15,7).Value = '''0.0021'
$hdr.Cells.Item(15,8).Value = 10
# row 16: 011244
$hdr.Cells.Item(16,1).Value = 14
$hdr.Cells.Item(16,2).Value = '''011244'
$hdr.Cells.Item(16,3).Value = '万家惠裕回报6个月持有期混合型证券投资基金C'
$hdr.Cells.Item(16,4).Value = '''0.14'
$hdr.Cells.Item(16,5).Value = '''23.04'
$hdr.Cells.Item(16,6).Value = '''0.88'
$hdr.Cells.Item(16,7).Value = '''0.0012'
$hdr.Cells.Item(16,8).Value = 6
# row 17: 005036
$hdr.Cells.Item(17,1).Value = 15
$hdr.Cells.Item(17,2).Value = '''005036'
$hdr.Cells.Item(17,3).Value = '银华信息科技量化优选股票C'
$hdr.Cells.Item(17,4).Value = '''0.04'
$hdr.Cells.Item(17,5).Value = '''90.98'
$hdr.Cells.Item(17,6).Value = '''2.28'
$hdr.Cells.Item(17,7).Value = '''0.0009'
$hdr.Cells.Item(17,8).Value = 6

# ============================ 总计 roll-up table ============================
$tot = $summarySheet
$styleSrc.Copy() | Out-Null
$tot.Range("B1:D1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$idxStyleSrc.Copy() | Out-Null
$tot.Range("A2:A7").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$tot.Cells.Item(1,2).Value = '日期'
$tot.Cells.Item(1,3).Value = '持有数量(只)'
$tot.Cells.Item(1,4).Value = '持有市值(亿元)'

$tot.Cells.Item(2,1).Value = 0
$tot.Cells.Item(2,2).Value = '2022-Q1'
$tot.Cells.Item(2,3).Value = 16
$tot.Cells.Item(2,4).Value = 1.17
$tot.Cells.Item(3,1).Value = 1
$tot.Cells.Item(3,2).Value = '2021-Q4'
$tot.Cells.Item(3,3).Value = 13
$tot.Cells.Item(3,4).Value = 1
$tot.Cells.Item(4,1).Value = 2
$tot.Cells.Item(4,2).Value = '2021-Q3'
$tot.Cells.Item(4,3).Value = 20
$tot.Cells.Item(4,4).Value = 1.88
$tot.Cells.Item(5,1).Value = 3
$tot.Cells.Item(5,2).Value = '2021-Q2'
$tot.Cells.Item(5,3).Value = 25
$tot.Cells.Item(5,4).Value = 8.96
$tot.Cells.Item(6,1).Value = 4
$tot.Cells.Item(6,2).Value = '2021-Q1'
$tot.Cells.Item(6,3).Value = 13
$tot.Cells.Item(6,4).Value = 5.1
$tot.Cells.Item(7,1).Value = 5
$tot.Cells.Item(7,2).Value = '2020-Q4'
$tot.Cells.Item(7,3).Value = 6
$tot.Cells.Item(7,4).Value = 0.51

# leave selection/view sane
$quarterSheet.Range("A1").Select() | Out-Null
